# Update 'F' column ('想去人数' / want-to-go count) values across all four
# worksheets to match the refreshed scrape output (gh-pages commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4517
$ws.Range("F5").Value = 3618
$ws.Range("F6").Value = 1041
$ws.Range("F7").Value = 165
$ws.Range("F9").Value = 355
$ws.Range("F10").Value = 352
$ws.Range("F11").Value = 2496
$ws.Range("F12").Value = 1276
$ws.Range("F14").Value = 1968
$ws.Range("F15").Value = 93
$ws.Range("F16").Value = 11
$ws.Range("F17").Value = 547
$ws.Range("F18").Value = 260
$ws.Range("F20").Value = 10326
$ws.Range("F21").Value = 6021
$ws.Range("F22").Value = 14
$ws.Range("F23").Value = 6
$ws.Range("F24").Value = 392
$ws.Range("F25").Value = 214
$ws.Range("F26").Value = 4
$ws.Range("F27").Value = 9
$ws.Range("F28").Value = 836
$ws.Range("F29").Value = 16
$ws.Range("F30").Value = 167
$ws.Range("F32").Value = 3557
$ws.Range("F36").Value = 120
$ws.Range("F37").Value = 258
$ws.Range("F39").Value = 239
$ws.Range("F40").Value = 4838
$ws.Range("F42").Value = 1118
$ws.Range("F43").Value = 162
$ws.Range("F44").Value = 131
$ws.Range("F45").Value = 90
$ws.Range("F46").Value = 483

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 27
$ws.Range("F15").Value = 3555

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8781
$ws.Range("F3").Value = 439
$ws.Range("F4").Value = 1617

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 439
$ws.Range("F3").Value = 1617
$ws.Range("F5").Value = 4517
$ws.Range("F8").Value = 3618
$ws.Range("F9").Value = 1041
$ws.Range("F10").Value = 165
$ws.Range("F12").Value = 352
$ws.Range("F13").Value = 2496
$ws.Range("F15").Value = 1277
$ws.Range("F18").Value = 97
$ws.Range("F19").Value = 11
$ws.Range("F21").Value = 547
$ws.Range("F22").Value = 260
$ws.Range("F24").Value = 10326
$ws.Range("F25").Value = 3555
$ws.Range("F27").Value = 14
$ws.Range("F28").Value = 392
$ws.Range("F29").Value = 214
$ws.Range("F30").Value = 4
$ws.Range("F31").Value = 9
$ws.Range("F32").Value = 836
$ws.Range("F33").Value = 16
$ws.Range("F34").Value = 167
$ws.Range("F36").Value = 3557
$ws.Range("F38").Value = 120
$ws.Range("F39").Value = 258
$ws.Range("F41").Value = 239
$ws.Range("F42").Value = 4838
$ws.Range("F44").Value = 1118
$ws.Range("F45").Value = 162
$ws.Range("F46").Value = 90
$ws.Range("F47").Value = 483

